# Apply change described in commit "#5: cash & deposit done":
# Extend the 存款 (deposit) worksheet with bank / deposit_type / currency
# headers (it previously reused a stray duplicate data row as a header) and
# append the common metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) used on the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) -------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Data rows (2-20): B-F already held the right data, just shifted down
# one row (the old row 1 was a duplicate of row 2 and is now the header), so
# rewrite B-F for completeness and fill in the new G-M metadata columns. ---

# Keep the new date column ("2011-11-22") stored as text, like the other
# sheets, instead of letting Excel auto-convert it to a date serial.
$ws.Range("I2:I20").NumberFormat = "@"

# row 2 (index 50)
$ws.Range("B2").Value = "台北富邦商業銀行松江分行"
$ws.Range("C2").Value = "活期存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "邱文彥"
$ws.Range("F2").Value = 2922
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2011-11-22"
$ws.Range("J2").Value = "邱文彥"
$ws.Range("K2").Value = 1743
$ws.Range("L2").Value = "tmpf3df1"
$ws.Range("M2").Value = 50

# row 3 (index 51)
$ws.Range("B3").Value = "基隆第一信用合作社八斗子分社"
$ws.Range("C3").Value = "活期存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "邱文彥"
$ws.Range("F3").Value = 519
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2011-11-22"
$ws.Range("J3").Value = "邱文彥"
$ws.Range("K3").Value = 1743
$ws.Range("L3").Value = "tmpf3df1"
$ws.Range("M3").Value = 51

# row 4 (index 52)
$ws.Range("B4").Value = "台新國際商業銀行古亭分行"
$ws.Range("C4").Value = "綜合存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "邱文彥"
$ws.Range("F4").Value = 1264
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2011-11-22"
$ws.Range("J4").Value = "邱文彥"
$ws.Range("K4").Value = 1743
$ws.Range("L4").Value = "tmpf3df1"
$ws.Range("M4").Value = 52

# row 5 (index 53)
$ws.Range("B5").Value = "灣銀行武昌分行"
$ws.Range("C5").Value = "活期存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "邱文彥"
$ws.Range("F5").Value = 226760
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2011-11-22"
$ws.Range("J5").Value = "邱文彥"
$ws.Range("K5").Value = 1743
$ws.Range("L5").Value = "tmpf3df1"
$ws.Range("M5").Value = 53

# row 6 (index 54)
$ws.Range("B6").Value = "臺灣銀行武昌分行"
$ws.Range("C6").Value = "活期存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "邱文彥"
$ws.Range("F6").Value = 84544
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2011-11-22"
$ws.Range("J6").Value = "邱文彥"
$ws.Range("K6").Value = 1743
$ws.Range("L6").Value = "tmpf3df1"
$ws.Range("M6").Value = 54

# row 7 (index 55)
$ws.Range("B7").Value = "臺灣銀行和平分行"
$ws.Range("C7").Value = "活期存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "邱文彥"
$ws.Range("F7").Value = 132
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2011-11-22"
$ws.Range("J7").Value = "邱文彥"
$ws.Range("K7").Value = 1743
$ws.Range("L7").Value = "tmpf3df1"
$ws.Range("M7").Value = 55

# row 8 (index 56)
$ws.Range("B8").Value = "f國信託商業銀行雙和分行"
$ws.Range("C8").Value = "活期存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "邱文彥"
$ws.Range("F8").Value = 8045
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2011-11-22"
$ws.Range("J8").Value = "邱文彥"
$ws.Range("K8").Value = 1743
$ws.Range("L8").Value = "tmpf3df1"
$ws.Range("M8").Value = 56

# row 9 (index 57)
$ws.Range("B9").Value = "中華郵政股份有限公司永和福和郵局"
$ws.Range("C9").Value = "活期存款"
$ws.Range("D9").Value = "新臺幣"
$ws.Range("E9").Value = "邱文彥"
$ws.Range("F9").Value = 108724
$ws.Range("G9").Value = "deposit"
$ws.Range("H9").Value = "normal"
$ws.Range("I9").Value = "2011-11-22"
$ws.Range("J9").Value = "邱文彥"
$ws.Range("K9").Value = 1743
$ws.Range("L9").Value = "tmpf3df1"
$ws.Range("M9").Value = 57

# row 10 (index 58)
$ws.Range("B10").Value = "中華郵政股份有限公司高雄西子灣郵局"
$ws.Range("C10").Value = "活期存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "邱文彥"
$ws.Range("F10").Value = 353603
$ws.Range("G10").Value = "deposit"
$ws.Range("H10").Value = "normal"
$ws.Range("I10").Value = "2011-11-22"
$ws.Range("J10").Value = "邱文彥"
$ws.Range("K10").Value = 1743
$ws.Range("L10").Value = "tmpf3df1"
$ws.Range("M10").Value = 58

# row 11 (index 60)
$ws.Range("B11").Value = "彰化商業銀行鹽埕分行"
$ws.Range("C11").Value = "活期存款"
$ws.Range("D11").Value = "新臺幣"
$ws.Range("E11").Value = "邱文彥"
$ws.Range("F11").Value = 70
$ws.Range("G11").Value = "deposit"
$ws.Range("H11").Value = "normal"
$ws.Range("I11").Value = "2011-11-22"
$ws.Range("J11").Value = "邱文彥"
$ws.Range("K11").Value = 1743
$ws.Range("L11").Value = "tmpf3df1"
$ws.Range("M11").Value = 60

# row 12 (index 61)
$ws.Range("B12").Value = "國泰世華商業銀行古亭分行"
$ws.Range("C12").Value = "活期存款"
$ws.Range("D12").Value = "新臺幣"
$ws.Range("E12").Value = "邱文彥"
$ws.Range("F12").Value = 145
$ws.Range("G12").Value = "deposit"
$ws.Range("H12").Value = "normal"
$ws.Range("I12").Value = "2011-11-22"
$ws.Range("J12").Value = "邱文彥"
$ws.Range("K12").Value = 1743
$ws.Range("L12").Value = "tmpf3df1"
$ws.Range("M12").Value = 61

# row 13 (index 62)
$ws.Range("B13").Value = "安泰商業銀行中和分行"
$ws.Range("C13").Value = "綜合存款"
$ws.Range("D13").Value = "新臺幣"
$ws.Range("E13").Value = "邱文彥"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "deposit"
$ws.Range("H13").Value = "normal"
$ws.Range("I13").Value = "2011-11-22"
$ws.Range("J13").Value = "邱文彥"
$ws.Range("K13").Value = 1743
$ws.Range("L13").Value = "tmpf3df1"
$ws.Range("M13").Value = 62

# row 14 (index 63)
$ws.Range("B14").Value = "安泰商業銀行中和分行"
$ws.Range("C14").Value = "其他存款"
$ws.Range("D14").Value = "新臺幣"
$ws.Range("E14").Value = "邱文彥"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "deposit"
$ws.Range("H14").Value = "normal"
$ws.Range("I14").Value = "2011-11-22"
$ws.Range("J14").Value = "邱文彥"
$ws.Range("K14").Value = 1743
$ws.Range("L14").Value = "tmpf3df1"
$ws.Range("M14").Value = 63

# row 15 (index 64)
$ws.Range("B15").Value = "合作金庫商業銀行大稻埕分行"
$ws.Range("C15").Value = "活期存款"
$ws.Range("D15").Value = "新臺幣"
$ws.Range("E15").Value = "黃淑芬"
$ws.Range("F15").Value = 37120
$ws.Range("G15").Value = "deposit"
$ws.Range("H15").Value = "normal"
$ws.Range("I15").Value = "2011-11-22"
$ws.Range("J15").Value = "邱文彥"
$ws.Range("K15").Value = 1743
$ws.Range("L15").Value = "tmpf3df1"
$ws.Range("M15").Value = 64

# row 16 (index 65)
$ws.Range("B16").Value = "臺灣中小企業銀行復興分行"
$ws.Range("C16").Value = "活期存款"
$ws.Range("D16").Value = "新臺幣"
$ws.Range("E16").Value = "黃淑芬"
$ws.Range("F16").Value = 103527
$ws.Range("G16").Value = "deposit"
$ws.Range("H16").Value = "normal"
$ws.Range("I16").Value = "2011-11-22"
$ws.Range("J16").Value = "邱文彥"
$ws.Range("K16").Value = 1743
$ws.Range("L16").Value = "tmpf3df1"
$ws.Range("M16").Value = 65

# row 17 (index 66)
$ws.Range("B17").Value = "臺灣中小企業銀行復興分行"
$ws.Range("C17").Value = "活期存款"
$ws.Range("D17").Value = "新臺幣"
$ws.Range("E17").Value = "黃淑芬"
$ws.Range("F17").Value = 512577
$ws.Range("G17").Value = "deposit"
$ws.Range("H17").Value = "normal"
$ws.Range("I17").Value = "2011-11-22"
$ws.Range("J17").Value = "邱文彥"
$ws.Range("K17").Value = 1743
$ws.Range("L17").Value = "tmpf3df1"
$ws.Range("M17").Value = 66

# row 18 (index 67)
$ws.Range("B18").Value = "臺灣中小企業銀行復興分行"
$ws.Range("C18").Value = "活期存款"
$ws.Range("D18").Value = "美金"
$ws.Range("E18").Value = "黃淑芬"
$ws.Range("F18").Value = 879.46
$ws.Range("G18").Value = "deposit"
$ws.Range("H18").Value = "normal"
$ws.Range("I18").Value = "2011-11-22"
$ws.Range("J18").Value = "邱文彥"
$ws.Range("K18").Value = 1743
$ws.Range("L18").Value = "tmpf3df1"
$ws.Range("M18").Value = 67

# row 19 (index 68)
$ws.Range("B19").Value = "台北富邦商業銀行安和分行"
$ws.Range("C19").Value = "活期存款"
$ws.Range("D19").Value = "新臺幣"
$ws.Range("E19").Value = "黃淑芬"
$ws.Range("F19").Value = 2025
$ws.Range("G19").Value = "deposit"
$ws.Range("H19").Value = "normal"
$ws.Range("I19").Value = "2011-11-22"
$ws.Range("J19").Value = "邱文彥"
$ws.Range("K19").Value = 1743
$ws.Range("L19").Value = "tmpf3df1"
$ws.Range("M19").Value = 68

# row 20 (index 69)
$ws.Range("B20").Value = "臺灣銀行"
$ws.Range("C20").Value = "其他存款"
$ws.Range("D20").Value = "新臺幣"
$ws.Range("E20").Value = "黃淑芬"
$ws.Range("F20").Value = 5287
$ws.Range("G20").Value = "deposit"
$ws.Range("H20").Value = "normal"
$ws.Range("I20").Value = "2011-11-22"
$ws.Range("J20").Value = "邱文彥"
$ws.Range("K20").Value = 1743
$ws.Range("L20").Value = "tmpf3df1"
$ws.Range("M20").Value = 69

# --- Styling: match the bold/bordered header style and the plain data-row
# style already used by columns B:F onto the newly added G:M columns. ------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("G2:M20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
